# Commit: "I0 and IF added"
# Adds two new columns (I = "I0", J = "IF") to the worksheet, extending
# the used range from A1:H69 to A1:J69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers, matching style of H1 ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell H1 onto the new
# header cells so they share the same cell style (bold, bordered,
# centered) as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2-69: numeric values for I0 and IF columns ---
$data = @(
    @{Row=2; I=9; J=9},
    @{Row=3; I=8; J=8},
    @{Row=4; I=9; J=9},
    @{Row=5; I=8; J=8},
    @{Row=6; I=9; J=9},
    @{Row=7; I=3; J=4},
    @{Row=8; I=5; J=6},
    @{Row=9; I=6; J=6},
    @{Row=10; I=7; J=7},
    @{Row=11; I=7; J=7},
    @{Row=12; I=6; J=7},
    @{Row=13; I=6; J=7},
    @{Row=14; I=7; J=7},
    @{Row=15; I=6; J=7},
    @{Row=16; I=6; J=7},
    @{Row=17; I=9; J=9},
    @{Row=18; I=9; J=9},
    @{Row=19; I=5; J=5},
    @{Row=20; I=6; J=6},
    @{Row=21; I=7; J=7},
    @{Row=22; I=9; J=9},
    @{Row=23; I=7; J=7},
    @{Row=24; I=9; J=9},
    @{Row=25; I=7; J=7},
    @{Row=26; I=7; J=7},
    @{Row=27; I=5; J=6},
    @{Row=28; I=7; J=7},
    @{Row=29; I=9; J=9},
    @{Row=30; I=8; J=8},
    @{Row=31; I=7; J=7},
    @{Row=32; I=7; J=7},
    @{Row=33; I=5; J=6},
    @{Row=34; I=4; J=5},
    @{Row=35; I=7; J=7},
    @{Row=36; I=6; J=6},
    @{Row=37; I=9; J=9},
    @{Row=38; I=9; J=9},
    @{Row=39; I=7; J=7},
    @{Row=40; I=9; J=9},
    @{Row=41; I=7; J=7},
    @{Row=42; I=8; J=8},
    @{Row=43; I=7; J=8},
    @{Row=44; I=6; J=7},
    @{Row=45; I=7; J=8},
    @{Row=46; I=6; J=7},
    @{Row=47; I=9; J=9},
    @{Row=48; I=8; J=8},
    @{Row=49; I=8; J=8},
    @{Row=50; I=6; J=6},
    @{Row=51; I=7; J=7},
    @{Row=52; I=7; J=8},
    @{Row=53; I=7; J=8},
    @{Row=54; I=7; J=7},
    @{Row=55; I=9; J=9},
    @{Row=56; I=7; J=7},
    @{Row=57; I=6; J=7},
    @{Row=58; I=7; J=7},
    @{Row=59; I=7; J=7},
    @{Row=60; I=7; J=8},
    @{Row=61; I=7; J=7},
    @{Row=62; I=8; J=8},
    @{Row=63; I=8; J=8},
    @{Row=64; I=7; J=8},
    @{Row=65; I=8; J=8},
    @{Row=66; I=7; J=7},
    @{Row=67; I=6; J=6},
    @{Row=68; I=8; J=8},
    @{Row=69; I=3; J=3}
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
}
